$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the latest scrape (GitHub Actions bot)
# Price column (D) values are textual (e.g. "1.300", "22.566.16") so force
# a Text number format before assignment to avoid Excel auto-converting them
# to numbers (which would normalize "19.80" -> 19.8, drop leading/trailing
# zeros, or use scientific notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.566.16"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.577.15"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.53"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3718"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.54"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3355"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.142"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07516"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.05"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.003"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.957"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.580.79"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.71"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06774"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.413"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.57"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.552.04"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.404"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.599"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.06"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.80"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.016"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.53"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.758.60"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.202"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.725"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08339"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02470"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.432"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06403"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.300"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.39"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6355"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.90"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6190"
$ws.Range("E46").Value = "  +5.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.793"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.068"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.43"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07284"
$ws.Range("E51").Value = "  -0.36%  "
